$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order line item appended as row 16 (SKU, Name, Quantity, Cost Per, Total Cost).
# Quantity/Cost/Total columns are stored as text in this sheet (like every other
# row), so force text formatting before writing the numeric-looking values to
# avoid Excel auto-converting them to the Number type.
$ws.Range("A16").Value = "TN380"
$ws.Range("B16").Value = "Natalie's - Strawberry Lemonade"

$ws.Range("C16:E16").NumberFormat = "@"
$ws.Range("C16").Value = "1"
$ws.Range("D16").Value = "10.15"
$ws.Range("E16").Value = "10.15"
$ws.Range("C16:E16").Style = "Normal"
